$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "request limit" style text values (shared strings)
$ws.Range("H7").Value = "-4:0"
$ws.Range("M7").Value = "3:0"

# Update the timestamp values (serial date/time) in the merged cells E4:F4 and J4:K4
$ws.Range("E4").Value = 45796.465090217825
$ws.Range("J4").Value = 45796.465090217825
